$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the new rows to be written as literal text (matches existing sheet convention
# where every cell - including numeric-looking SKUs/quantities/prices - is stored as text).
$dataRange = $ws.Range("A3:E17")
$dataRange.NumberFormat = "@"

$ws.Range("A3").Value = "711SPRNKLEBL"
$ws.Range("B3").Value = "Sprinkles - Blue"
$ws.Range("C3").Value = "2"
$ws.Range("D3").Value = "24.99"
$ws.Range("E3").Value = "49.98"

$ws.Range("A4").Value = "711SPRKLCH25"
$ws.Range("B4").Value = "Sprinkles - Dark Chocolate"
$ws.Range("C4").Value = "2"
$ws.Range("D4").Value = "59.49"
$ws.Range("E4").Value = "118.98"

$ws.Range("A5").Value = "711SPRNKLEOR"
$ws.Range("B5").Value = "Sprinkles - Orange"
$ws.Range("C5").Value = "1"
$ws.Range("D5").Value = "23.67"
$ws.Range("E5").Value = "23.67"

$ws.Range("A6").Value = "10207579"
$ws.Range("B6").Value = "Salt - Sea Coarse"
$ws.Range("C6").Value = "4"
$ws.Range("D6").Value = "18.72"
$ws.Range("E6").Value = "74.88"

$ws.Range("A7").Value = "24510105CB"
$ws.Range("B7").Value = "Box Cake - 10x10x5"
$ws.Range("C7").Value = "1"
$ws.Range("D7").Value = "59.49"
$ws.Range("E7").Value = "59.49"

$ws.Range("A8").Value = "245885CB"
$ws.Range("B8").Value = "Box Cake - 8x8x5"
$ws.Range("C8").Value = "1"
$ws.Range("D8").Value = "34.81"
$ws.Range("E8").Value = "34.81"

$ws.Range("A9").Value = "770V9I30008"
$ws.Range("B9").Value = "Eclair Paper"
$ws.Range("C9").Value = "3"
$ws.Range("D9").Value = "88.99"
$ws.Range("E9").Value = "266.97"

$ws.Range("A10").Value = "245CCGR1410BL"
$ws.Range("B10").Value = "Cake Board - 1/4 Sheet"
$ws.Range("C10").Value = "2"
$ws.Range("D10").Value = "40.49"
$ws.Range("E10").Value = "80.98"

$ws.Range("A11").Value = "24510102WB"
$ws.Range("B11").Value = "Box Cake - 10x10x2.5 Window"
$ws.Range("C11").Value = "4"
$ws.Range("D11").Value = "65.47"
$ws.Range("E11").Value = "261.88"

$ws.Range("A12").Value = "245882WB"
$ws.Range("B12").Value = "Box Cake - 8x8x2.5 (window)"
$ws.Range("C12").Value = "2"
$ws.Range("D12").Value = "68.60"
$ws.Range("E12").Value = "137.20"

$ws.Range("A13").Value = "150300865"
$ws.Range("B13").Value = "Bag Paper - 6x13.5 Window"
$ws.Range("C13").Value = "2"
$ws.Range("D13").Value = "79.99"
$ws.Range("E13").Value = "159.98"

$ws.Range("A14").Value = "150BB6218N"
$ws.Range("B14").Value = "Bag Paper - 6.5x17.75 Window"
$ws.Range("C14").Value = "3"
$ws.Range("D14").Value = "104.99"
$ws.Range("E14").Value = "314.97"

$ws.Range("A15").Value = "130PPF0628M1M"
$ws.Range("B15").Value = "Bag - Perforated (6x28)"
$ws.Range("C15").Value = "1"
$ws.Range("D15").Value = "54.99"
$ws.Range("E15").Value = "54.99"

$ws.Range("A16").Value = "433qlinerbl"
$ws.Range("B16").Value = "Sheet Pan Liner - White"
$ws.Range("C16").Value = "3"
$ws.Range("D16").Value = "43.99"
$ws.Range("E16").Value = "131.97"

$ws.Range("A17").Value = "656095131"
$ws.Range("B17").Value = "Container - Muffin (12 Pack)"
$ws.Range("C17").Value = "2"
$ws.Range("D17").Value = "65.49"
$ws.Range("E17").Value = "130.98"

# Reset to the default cell style so no explicit style index is left on the new cells
# (keeps the written XML free of spurious s="..." attributes), while the values remain text.
$dataRange.Style = "Normal"

$ws.Range("A1").Select()